$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2157772621809745
$ws.Range("C2").Value = 0.5197215777262181
$ws.Range("J2").Value = 0.0185614849187935
$ws.Range("P2").Value = 0.1577726218097448
$ws.Range("S2").Value = 0.08816705336426914
$ws.Range("B3").Value = 0.01304347826086956
$ws.Range("C3").Value = 0.01739130434782609
$ws.Range("J3").Value = 0.02173913043478261
$ws.Range("P3").Value = 0.7043478260869566
$ws.Range("S3").Value = 0.2434782608695652
$ws.Range("J4").Value = 0.07407407407407407
$ws.Range("P4").Value = 0.7222222222222222
$ws.Range("S4").Value = 0.2037037037037037
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.05118110236220472
$ws.Range("D6").Value = 0.01181102362204724
$ws.Range("F6").Value = 0.03543307086614173
$ws.Range("J6").Value = 0.3031496062992126
$ws.Range("Q6").Value = 0.1220472440944882
$ws.Range("R6").Value = 0.09055118110236221
$ws.Range("S6").Value = 0.3858267716535433
$ws.Range("B7").Value = 0.1345454545454546
$ws.Range("D7").Value = 0.02545454545454546
$ws.Range("F7").Value = 0.04727272727272727
$ws.Range("J7").Value = 0.1018181818181818
$ws.Range("O7").Value = 0.01818181818181818
$ws.Range("Q7").Value = 0.2036363636363636
$ws.Range("R7").Value = 0.08363636363636363
$ws.Range("S7").Value = 0.3854545454545454
$ws.Range("B8").Value = 0.1054945054945055
$ws.Range("D8").Value = 0.01758241758241758
$ws.Range("E8").Value = 0.002197802197802198
$ws.Range("F8").Value = 0.06153846153846154
$ws.Range("J8").Value = 0.1296703296703297
$ws.Range("O8").Value = 0.02417582417582418
$ws.Range("Q8").Value = 0.1824175824175824
$ws.Range("R8").Value = 0.09230769230769231
$ws.Range("S8").Value = 0.3846153846153846
$ws.Range("B9").Value = 0.07053941908713693
$ws.Range("D9").Value = 0.02489626556016597
$ws.Range("F9").Value = 0.05809128630705394
$ws.Range("J9").Value = 0.1286307053941909
$ws.Range("O9").Value = 0.03319502074688797
$ws.Range("Q9").Value = 0.2033195020746888
$ws.Range("R9").Value = 0.0912863070539419
$ws.Range("S9").Value = 0.3900414937759336
$ws.Range("B10").Value = 0.1413404114134041
$ws.Range("D10").Value = 0.019907100199071
$ws.Range("E10").Value = 0.0006635700066357001
$ws.Range("F10").Value = 0.0650298606502986
$ws.Range("J10").Value = 0.1220968812209688
$ws.Range("O10").Value = 0.0232249502322495
$ws.Range("Q10").Value = 0.2289316522893165
$ws.Range("R10").Value = 0.0583941605839416
$ws.Range("S10").Value = 0.3404114134041141
$ws.Range("G11").Value = 0.1666666666666667
$ws.Range("J11").Value = 0.06862745098039216
$ws.Range("K11").Value = 0.1838235294117647
$ws.Range("L11").Value = 0.5612745098039216
$ws.Range("S11").Value = 0.0196078431372549
$ws.Range("G12").Value = 0.7939914163090128
$ws.Range("J12").Value = 0.1587982832618026
$ws.Range("K12").Value = 0.0128755364806867
$ws.Range("L12").Value = 0.01716738197424893
$ws.Range("S12").Value = 0.01716738197424893
$ws.Range("G13").Value = 0.5576923076923077
$ws.Range("J13").Value = 0.3461538461538461
$ws.Range("S13").Value = 0.09615384615384616
$ws.Range("F15").Value = 0.02008032128514056
$ws.Range("H15").Value = 0.1485943775100401
$ws.Range("I15").Value = 0.07630522088353414
$ws.Range("J15").Value = 0.3253012048192771
$ws.Range("K15").Value = 0.06425702811244979
$ws.Range("M15").Value = 0.008032128514056224
$ws.Range("N15").Value = 0.004016064257028112
$ws.Range("O15").Value = 0.04819277108433735
$ws.Range("S15").Value = 0.3052208835341366
$ws.Range("F16").Value = 0.04182509505703422
$ws.Range("H16").Value = 0.1368821292775665
$ws.Range("I16").Value = 0.07604562737642585
$ws.Range("J16").Value = 0.4448669201520912
$ws.Range("K16").Value = 0.1330798479087452
$ws.Range("M16").Value = 0.01140684410646388
$ws.Range("O16").Value = 0.04562737642585551
$ws.Range("S16").Value = 0.1102661596958175
$ws.Range("F17").Value = 0.02491103202846975
$ws.Range("H17").Value = 0.1405693950177936
$ws.Range("I17").Value = 0.1192170818505338
$ws.Range("J17").Value = 0.4181494661921708
$ws.Range("K17").Value = 0.1067615658362989
$ws.Range("M17").Value = 0.01423487544483986
$ws.Range("O17").Value = 0.06761565836298933
$ws.Range("S17").Value = 0.1085409252669039
$ws.Range("F18").Value = 0.03076923076923077
$ws.Range("H18").Value = 0.1384615384615385
$ws.Range("I18").Value = 0.05641025641025641
$ws.Range("J18").Value = 0.4461538461538462
$ws.Range("K18").Value = 0.1179487179487179
$ws.Range("M18").Value = 0.02051282051282051
$ws.Range("O18").Value = 0.07692307692307693
$ws.Range("S18").Value = 0.1128205128205128
$ws.Range("F19").Value = 0.01780821917808219
$ws.Range("H19").Value = 0.1972602739726027
$ws.Range("I19").Value = 0.08424657534246575
$ws.Range("J19").Value = 0.360958904109589
$ws.Range("K19").Value = 0.126027397260274
$ws.Range("M19").Value = 0.02602739726027397
$ws.Range("N19").Value = 0.0006849315068493151
$ws.Range("O19").Value = 0.0589041095890411
$ws.Range("S19").Value = 0.1280821917808219
